$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, `
                             $true, 1, $false, $new, 2)
}

Replace-Text "2025-04-09 Wednesday" "2025-04-10 Thursday"

Replace-Text "865÷3=" "926÷8="
Replace-Text "237÷4=" "557÷8="
Replace-Text "195÷5=" "277÷9="
Replace-Text "264÷4=" "167÷9="
Replace-Text "655÷9=" "327÷8="
Replace-Text "306÷8=" "248÷8="
Replace-Text "691÷8=" "658÷9="
Replace-Text "801÷4=" "400÷4="
Replace-Text "659÷5=" "118÷4="
Replace-Text "114÷4=" "230÷2="
Replace-Text "687÷3=" "153÷2="
Replace-Text "638÷2=" "785÷3="
Replace-Text "109÷6=" "807÷9="
Replace-Text "918÷9=" "677÷4="
Replace-Text "402÷6=" "702÷5="
Replace-Text "976÷9=" "394÷7="
Replace-Text "929÷9=" "181÷4="
Replace-Text "446÷9=" "587÷7="
Replace-Text "766÷8=" "120÷5="
Replace-Text "281÷5=" "937÷9="
Replace-Text "651÷4=" "194÷3="
Replace-Text "866÷7=" "130÷5="
Replace-Text "577÷5=" "592÷4="
Replace-Text "973÷8=" "691÷9="
Replace-Text "897÷3=" "327÷7="
